$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 273; all existing rows 273-326 shift down to 274-327.
$ws.Rows(273).Insert()

$ws.Range("A273").Value = 3
$ws.Range("B273").Value = "Femacal de La Calera"
$ws.Range("C273").Value = "Coquimbo"
$ws.Range("D273").Value = (Get-Date -Year 2021 -Month 10 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E273").Value = 5
$ws.Range("F273").Value = 100112045
$ws.Range("G273").Value = "Zapallo"
$ws.Range("H273").Value = "Camote"
$ws.Range("I273").Value = "1a (guarda)"
$ws.Range("J273").Value = 200
$ws.Range("K273").Value = 650
$ws.Range("L273").Value = 700
$ws.Range("M273").Value = 678
$ws.Range("N273").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O273").Value = "Provincia de Talca"
$ws.Range("P273").Value = 678
$ws.Range("Q273").Value = 1
$ws.Range("R273").Value = "Hortaliza"
